$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is an empty placeholder row with no data; clearing its (nonexistent)
# contents drops the now-pointless empty <row r="2"/> element from the sheet
# without shifting any of the data in the rows below it.
$ws.Rows.Item(2).ClearContents()

# Row 4 ("nilesh") has now been settled - update its figures accordingly.
$ws.Range("B4").Value = 138000
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 12000
$ws.Range("G4").Value = 30000
$ws.Range("H4").Value = 12000
